$wb = $excel.ActiveWorkbook

# Row 17 (One for the Road) - ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1684802.2
$ws.Range("J17").Value = 1749552.4
$ws.Range("L17").Value = 5248657.199999999
$ws.Range("N17").Value = -5248993.199999999

# Row 32 (Automata for the People) - ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1403.4706
$ws.Range("I32").Value = 1158
$ws.Range("J32").Value = 1537.3636
$ws.Range("K32").Value = 1158
$ws.Range("L32").Value = 1537.3636
$ws.Range("M32").Value = -832
$ws.Range("N32").Value = -2189.3636

# Row 112 (Making Ends Meet) - ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2623.7222
$ws.Range("J112").Value = 2242.3125
$ws.Range("L112").Value = 6726.9375
$ws.Range("N112").Value = -8942.9375

# Row 138 (All-night Crafting) - ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 426939.38
$ws.Range("I138").Value = 732360.4
$ws.Range("J138").Value = 4048.7693
$ws.Range("K138").Value = 2197081.2
$ws.Range("L138").Value = 12146.3079
$ws.Range("M138").Value = -2191941.2
$ws.Range("N138").Value = -22426.3079

# Row 2 (Ain't Got No Ingots) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3739.45
$ws.Range("J2").Value = 3337.5
$ws.Range("L2").Value = 3337.5
$ws.Range("N2").Value = -3563.5

# Row 45 (Hollow Hallmarks) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 85026.48
$ws.Range("I45").Value = 129491.375
$ws.Range("K45").Value = 129491.375
$ws.Range("M45").Value = -129114.375

# Row 61 (Dealing with the Tough Stuff) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 26124.125
$ws.Range("I61").Value = 58332.332
$ws.Range("K61").Value = 58332.332
$ws.Range("M61").Value = -58120.332

# Row 74 (As the Bolt Flies) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9584.467000000001
$ws.Range("I74").Value = 13274.111
$ws.Range("K74").Value = 13274.111
$ws.Range("M74").Value = -12400.111

# Row 77 (Heavy Metal Banned (L)) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9584.467000000001
$ws.Range("I77").Value = 13274.111
$ws.Range("K77").Value = 66370.55500000001
$ws.Range("M77").Value = -62002.55500000001

# Row 94 (Setting the Stage) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 264001000
$ws.Range("J94").Value = 264001000
$ws.Range("L94").Value = 264001000
$ws.Range("N94").Value = -264002802

# Row 102 (Smells of Rich Tama-hagane) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 8326.666999999999
$ws.Range("I102").Value = 10261.375
$ws.Range("K102").Value = 10261.375
$ws.Range("M102").Value = -8639.375

# Row 116 (No Scope) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3739.45
$ws.Range("J116").Value = 3337.5
$ws.Range("L116").Value = 3337.5
$ws.Range("N116").Value = -7925.5

# Row 132 (Don't Bore Me, Ore Me) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4564.4443
$ws.Range("I132").Value = 4546.2856
$ws.Range("J132").Value = 4628
$ws.Range("K132").Value = 13638.8568
$ws.Range("L132").Value = 13884
$ws.Range("M132").Value = -11108.8568
$ws.Range("N132").Value = -18944

# Row 134 (Brace for More Vambraces) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 400000
$ws.Range("J134").Value = 400000
$ws.Range("L134").Value = 400000
$ws.Range("N134").Value = -410140

# Row 135 (Forgiveness for My Shins) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 208310.17
$ws.Range("J135").Value = 208310.17
$ws.Range("L135").Value = 208310.17
$ws.Range("N135").Value = -218450.17

# Row 136 (Metal with Mettle) - ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 26124.125
$ws.Range("I136").Value = 58332.332
$ws.Range("K136").Value = 174996.996
$ws.Range("M136").Value = -172446.996

# Row 3 (Hells Bells) - BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3739.45
$ws.Range("J3").Value = 3337.5
$ws.Range("L3").Value = 3337.5
$ws.Range("N3").Value = -3565.5

# Row 20 (Smelt and Dealt) - BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2403.4211
$ws.Range("I20").Value = 1578.6
$ws.Range("J20").Value = 3319.889
$ws.Range("K20").Value = 1578.6
$ws.Range("L20").Value = 3319.889
$ws.Range("M20").Value = -1331.6
$ws.Range("N20").Value = -3813.889

# Row 94 (High Steal) - BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8882.388999999999
$ws.Range("I94").Value = 11102.77
$ws.Range("J94").Value = 3109.4
$ws.Range("K94").Value = 11102.77
$ws.Range("L94").Value = 3109.4
$ws.Range("M94").Value = -10651.77
$ws.Range("N94").Value = -4011.4

# Row 99 (Meddle in Metal) - BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 10207.344
$ws.Range("I99").Value = 10662.654
$ws.Range("K99").Value = 10662.654
$ws.Range("M99").Value = -9164.654

# Row 110 (Selective Logging) - BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180

# Row 112 (Enlistment Highs) - BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 134 (Ruthenium Supremium) - BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 12980.454
$ws.Range("I134").Value = 12980.454
$ws.Range("K134").Value = 38941.362
$ws.Range("M134").Value = -36406.362

# Row 31 (Wall Not Found) - CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10020.611
$ws.Range("I31").Value = 11531.583
$ws.Range("J31").Value = 6998.6665
$ws.Range("K31").Value = 11531.583
$ws.Range("L31").Value = 6998.6665
$ws.Range("M31").Value = -11236.583
$ws.Range("N31").Value = -7588.6665

# Row 34 (Armoires of the Rich and Famous) - CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10020.611
$ws.Range("I34").Value = 11531.583
$ws.Range("J34").Value = 6998.6665
$ws.Range("K34").Value = 11531.583
$ws.Range("L34").Value = 6998.6665
$ws.Range("M34").Value = -11329.583
$ws.Range("N34").Value = -7402.6665

# Row 86 (Birch, Please) - CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 12001.6
$ws.Range("I86").Value = 9500
$ws.Range("K86").Value = 9500
$ws.Range("M86").Value = -8377

# Row 89 (Built This City on Blocks and Soul (L)) - CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 12001.6
$ws.Range("I89").Value = 9500
$ws.Range("K89").Value = 47500
$ws.Range("M89").Value = -41884

# Row 109 (Playing the Market) - CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 49991.5
$ws.Range("J109").Value = 49991.5
$ws.Range("L109").Value = 49991.5
$ws.Range("N109").Value = -52071.5

# Row 132 (Hull Lotta Damage) - CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1246.6578
$ws.Range("I132").Value = 1193.5428
$ws.Range("K132").Value = 3580.6284
$ws.Range("M132").Value = -1050.6284

# Row 134 (Wood You Be Quiet) - CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6950.5
$ws.Range("I134").Value = 8013.1875
$ws.Range("J134").Value = 2699.75
$ws.Range("K134").Value = 24039.5625
$ws.Range("L134").Value = 8099.25
$ws.Range("M134").Value = -21504.5625
$ws.Range("N134").Value = -13169.25

# Row 2 (Pork Is a Salty Food) - CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 81.333336
$ws.Range("I2").Value = 80.14286
$ws.Range("K2").Value = 480.85716
$ws.Range("M2").Value = -367.85716

# Row 18 (Fisher of Men) - CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 4273.25
$ws.Range("I18").Value = 5484.8335
$ws.Range("J18").Value = 638.5
$ws.Range("K18").Value = 16454.5005
$ws.Range("L18").Value = 1915.5
$ws.Range("M18").Value = -16285.5005
$ws.Range("N18").Value = -2253.5

# Row 23 (Sweet Smell of Success) - CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 23810104
$ws.Range("J23").Value = 27778322
$ws.Range("L23").Value = 83334966
$ws.Range("N23").Value = -83335436

# Row 33 (Cooking with Gas) - CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 259.85715
$ws.Range("J33").Value = 254.2
$ws.Range("L33").Value = 1525.2
$ws.Range("N33").Value = -2091.2

# Row 139 (Najoothie) - CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1113809.8
$ws.Range("I139").Value = 2308989.5
$ws.Range("K139").Value = 6926968.5
$ws.Range("M139").Value = -6921828.5

# Row 140 (Sweet, Sweet Bean Juice) - CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 13826.357
$ws.Range("I140").Value = 14524.538
$ws.Range("K140").Value = 43573.614
$ws.Range("M140").Value = -38393.614

# Row 141 (Ocean Explosion) - CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2717.1428
$ws.Range("I141").Value = 2586.6667
$ws.Range("K141").Value = 7760.000100000001
$ws.Range("M141").Value = -2580.000100000001

# Row 7 (Tan Before the Ban) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18010.775
$ws.Range("I7").Value = 22293.932
$ws.Range("J7").Value = 6718.8184
$ws.Range("K7").Value = 22293.932
$ws.Range("L7").Value = 6718.8184
$ws.Range("M7").Value = -22181.932
$ws.Range("N7").Value = -6942.8184

# Row 22 (Skin off Their Backs) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12976.471
$ws.Range("I22").Value = 22688.889
$ws.Range("K22").Value = 22688.889
$ws.Range("M22").Value = -22393.889

# Row 27 (Fire and Hide) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 12976.471
$ws.Range("I27").Value = 22688.889
$ws.Range("K27").Value = 22688.889
$ws.Range("M27").Value = -22581.889

# Row 40 (Best Served Toad) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 20109.371
$ws.Range("I40").Value = 24617.15
$ws.Range("J40").Value = 14099
$ws.Range("K40").Value = 24617.15
$ws.Range("L40").Value = 14099
$ws.Range("M40").Value = -24481.15
$ws.Range("N40").Value = -14371

# Row 93 (Hide to Go Seek) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3552.6086
$ws.Range("I93").Value = 3985.8948
$ws.Range("K93").Value = 3985.8948
$ws.Range("M93").Value = -2737.8948

# Row 122 (Hell on Leather) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5832.885
$ws.Range("I122").Value = 5819.8823
$ws.Range("J122").Value = 5857.4443
$ws.Range("K122").Value = 17459.6469
$ws.Range("L122").Value = 17572.3329
$ws.Range("M122").Value = -15009.6469
$ws.Range("N122").Value = -22472.3329

# Row 126 (Battered Books) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 18010.775
$ws.Range("I126").Value = 22293.932
$ws.Range("J126").Value = 6718.8184
$ws.Range("K126").Value = 66881.796
$ws.Range("L126").Value = 20156.4552
$ws.Range("M126").Value = -64411.796
$ws.Range("N126").Value = -25096.4552

# Row 135 (Dreams of Ja) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 112192.336
$ws.Range("J135").Value = 112192.336
$ws.Range("L135").Value = 112192.336
$ws.Range("N135").Value = -122332.336

# Row 136 (Respect for Br'aax) - LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4912.077
$ws.Range("I136").Value = 4210.7
$ws.Range("K136").Value = 12632.1
$ws.Range("M136").Value = -10082.1

# Row 32 (Piling It On) - WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10000
$ws.Range("I32").Value = 10000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9683
$ws.Range("N32").ClearContents()

# Row 96 (Skills on Display) - WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1730
$ws.Range("J96").Value = 4500
$ws.Range("L96").Value = 4500
$ws.Range("N96").Value = -7246

# Row 122 (Heavy Armoire) - WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3972.5435
$ws.Range("I122").Value = 2237.8965
$ws.Range("K122").Value = 6713.689499999999
$ws.Range("M122").Value = -4263.689499999999

# Row 126 (A Polished Purchase) - WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 18723.625
$ws.Range("I126").Value = 23916.055
$ws.Range("J126").Value = 3146.3333
$ws.Range("K126").Value = 71748.16500000001
$ws.Range("L126").Value = 9438.999899999999
$ws.Range("M126").Value = -69278.16500000001
$ws.Range("N126").Value = -14378.9999
